$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 272, pushing existing rows 272-277 down to 274-279.
$ws.Rows("272:273").Insert()

# New row 272: Red Globe bandeja data
$ws.Range("A272").Value = 4
$ws.Range("B272").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C272").Value = "Los Lagos"
$ws.Range("D272").Value = 44939
$ws.Range("E272").Value = 10
$ws.Range("F272").Value = "Fruta"
$ws.Range("G272").Value = 100109
$ws.Range("H272").Value = "Uva"
$ws.Range("I272").Value = 100109001
$ws.Range("J272").Value = "Uva"
$ws.Range("K272").Value = "Red Globe"
$ws.Range("L272").Value = "Primera"
$ws.Range("M272").Value = 400
$ws.Range("N272").Value = 16000
$ws.Range("O272").Value = 17000
$ws.Range("P272").Value = 16500
$ws.Range("Q272").Value = "`$/bandeja 10 kilos"
$ws.Range("R272").Value = "Provincia de Limarí"
$ws.Range("S272").Value = 1650
$ws.Range("T272").Value = 10

# New row 273: Superior Seedless bandeja data
$ws.Range("A273").Value = 4
$ws.Range("B273").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C273").Value = "Los Lagos"
$ws.Range("D273").Value = 44939
$ws.Range("E273").Value = 10
$ws.Range("F273").Value = "Fruta"
$ws.Range("G273").Value = 100109
$ws.Range("H273").Value = "Uva"
$ws.Range("I273").Value = 100109001
$ws.Range("J273").Value = "Uva"
$ws.Range("K273").Value = "Superior Seedless"
$ws.Range("L273").Value = "Primera"
$ws.Range("M273").Value = 400
$ws.Range("N273").Value = 15000
$ws.Range("O273").Value = 16000
$ws.Range("P273").Value = 15500
$ws.Range("Q273").Value = "`$/bandeja 10 kilos"
$ws.Range("R273").Value = "Provincia de Limarí"
$ws.Range("S273").Value = 1550
$ws.Range("T273").Value = 10

# Apply same date style (numFmt) as the rest of column D to the new D cells.
$ws.Range("D272").NumberFormat = $ws.Range("D274").NumberFormat
$ws.Range("D273").NumberFormat = $ws.Range("D274").NumberFormat

$wb.Save()
